$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 10.12531930373934
$ws.Range("C2").Value = 4.673205954955797
$ws.Range("E2").Value = 12.14585239965294
$ws.Range("F2").Value = 16.86991607391233
$ws.Range("G2").Value = 35.17305370159212
$ws.Range("H2").Value = 15.82555108395776
$ws.Range("K2").Value = 9.455566637351513
$ws.Range("M2").Value = 14.08704352671802
$ws.Range("N2").Value = 19.60505950283795

$ws.Range("B3").Value = 9.879514157090451
$ws.Range("C3").Value = 4.496068194755836
$ws.Range("E3").Value = 11.92214965718733
$ws.Range("F3").Value = 15.89584955866808
$ws.Range("G3").Value = 35.15954875482266
$ws.Range("H3").Value = 15.86928600165002
$ws.Range("K3").Value = 9.297095232072641
$ws.Range("M3").Value = 13.93117522771406
$ws.Range("N3").Value = 19.67257200615339

$ws.Range("B4").Value = 9.728095716980095
$ws.Range("C4").Value = 4.382668462738046
$ws.Range("E4").Value = 11.78671725173012
$ws.Range("F4").Value = 15.26997757108491
$ws.Range("G4").Value = 35.16259127605312
$ws.Range("H4").Value = 15.89896148038926
$ws.Range("K4").Value = 9.200719138236048
$ws.Range("M4").Value = 13.83834632795204
$ws.Range("N4").Value = 19.71595641163805

$ws.Range("B5").Value = 9.666367563570274
$ws.Range("C5").Value = 4.335328207892892
$ws.Range("E5").Value = 11.73209192024214
$ws.Range("F5").Value = 15.00819731993403
$ws.Range("G5").Value = 35.16667514078919
$ws.Range("H5").Value = 15.91176285084509
$ws.Range("K5").Value = 9.16173167569981
$ws.Range("M5").Value = 13.80128308892574
$ws.Range("N5").Value = 19.73412289204618

$ws.Range("B6").Value = 9.656119198528204
$ws.Range("C6").Value = 4.327400382210422
$ws.Range("E6").Value = 11.72305788199433
$ws.Range("F6").Value = 14.96433081551593
$ws.Range("G6").Value = 35.16752480606207
$ws.Range("H6").Value = 15.91393125509798
$ws.Range("K6").Value = 9.155276715756102
$ws.Range("M6").Value = 13.79517619427763
$ws.Range("N6").Value = 19.73716887321605

$ws.Range("B7").Value = 9.727263191551408
$ws.Range("C7").Value = 4.38203453425351
$ws.Range("E7").Value = 11.78597816291827
$ws.Range("F7").Value = 15.26647399323137
$ws.Range("G7").Value = 35.16263484725886
$ws.Range("H7").Value = 15.89913125769687
$ws.Range("K7").Value = 9.200192108362289
$ws.Range("M7").Value = 13.83784332674655
$ws.Range("N7").Value = 19.71619943766469

$ws.Range("B8").Value = 10.04072694349043
$ws.Range("C8").Value = 4.613111848515545
$ws.Range("E8").Value = 12.06836955305944
$ws.Range("F8").Value = 16.53996406344768
$ws.Range("G8").Value = 35.16604172717707
$ws.Range("H8").Value = 15.84004441808055
$ws.Range("K8").Value = 9.400765710411305
$ws.Range("M8").Value = 14.03272936669931
$ws.Range("N8").Value = 19.62793772834851

$ws.Range("B9").Value = 10.64731802559112
$ws.Range("C9").Value = 5.028075904327005
$ws.Range("E9").Value = 12.63382639572961
$ws.Range("F9").Value = 19.00274580682531
$ws.Range("G9").Value = 35.26282141128744
$ws.Range("H9").Value = 15.74662252411837
$ws.Range("K9").Value = 9.799090276058118
$ws.Range("M9").Value = 14.43577777512291
$ws.Range("N9").Value = 19.47012132967644

$ws.Range("B10").Value = 11.0827454403407
$ws.Range("C10").Value = 5.308108242973054
$ws.Range("E10").Value = 13.05169228478083
$ws.Range("F10").Value = 20.67494806633232
$ws.Range("G10").Value = 35.38891218123201
$ws.Range("H10").Value = 15.69174437983032
$ws.Range("K10").Value = 10.09169757178276
$ws.Range("M10").Value = 14.74196086304866
$ws.Range("N10").Value = 19.36339260122981

$ws.Range("B11").Value = 11.27759745017972
$ws.Range("C11").Value = 5.42985028304418
$ws.Range("E11").Value = 13.24135920500874
$ws.Range("F11").Value = 21.3917225636224
$ws.Range("G11").Value = 35.45817146880849
$ws.Range("H11").Value = 15.66978173301484
$ws.Range("K11").Value = 10.22417375007419
$ws.Range("M11").Value = 14.88285845109118
$ws.Range("N11").Value = 19.3168222979076

$ws.Range("B12").Value = 11.3508424083319
$ws.Range("C12").Value = 5.47512145836225
$ws.Range("E12").Value = 13.31304441543563
$ws.Range("F12").Value = 21.65686569030329
$ws.Range("G12").Value = 35.48610126827212
$ws.Range("H12").Value = 15.66189787708267
$ws.Range("K12").Value = 10.27419802769703
$ws.Range("M12").Value = 14.93639619425703
$ws.Range("N12").Value = 19.29947087901255

$ws.Range("B13").Value = 11.33509303745339
$ws.Range("C13").Value = 5.465408668020431
$ws.Range("E13").Value = 13.29761302525369
$ws.Range("F13").Value = 21.60004134736742
$ws.Range("G13").Value = 35.48001051567191
$ws.Range("H13").Value = 15.66357653235223
$ws.Range("K13").Value = 10.26343148428216
$ws.Range("M13").Value = 14.9248585270258
$ws.Range("N13").Value = 19.303195215909

$ws.Range("B14").Value = 11.28363462669857
$ws.Range("C14").Value = 5.433591495526384
$ws.Range("E14").Value = 13.24725997431332
$ws.Range("F14").Value = 21.4136618050453
$ws.Range("G14").Value = 35.46043520105657
$ws.Range("H14").Value = 15.66912443989305
$ws.Range("K14").Value = 10.22829239725361
$ws.Range("M14").Value = 14.88725964801583
$ws.Range("N14").Value = 19.31538910605136

$ws.Range("B15").Value = 11.25204218101047
$ws.Range("C15").Value = 5.413993994931379
$ws.Range("E15").Value = 13.21639714554385
$ws.Range("F15").Value = 21.29868154950795
$ws.Range("G15").Value = 35.44866621650817
$ws.Range("H15").Value = 15.67257910978456
$ws.Range("K15").Value = 10.2067487883451
$ws.Range("M15").Value = 14.86425161511957
$ws.Range("N15").Value = 19.32289513415849

$ws.Range("B16").Value = 11.06993973672541
$ws.Range("C16").Value = 5.300036899060774
$ws.Range("E16").Value = 13.03928196302569
$ws.Range("F16").Value = 20.62722412089977
$ws.Range("G16").Value = 35.38462472644818
$ws.Range("H16").Value = 15.69324020298165
$ws.Range("K16").Value = 10.08302283947189
$ws.Range("M16").Value = 14.73278113261737
$ws.Range("N16").Value = 19.36647584769706

$ws.Range("B17").Value = 10.95734256973027
$ws.Range("C17").Value = 5.228668118442291
$ws.Range("E17").Value = 12.93046363960222
$ws.Range("F17").Value = 20.20408069597325
$ws.Range("G17").Value = 35.34838010733903
$ws.Range("H17").Value = 15.70668484906052
$ws.Range("K17").Value = 10.00692231412899
$ws.Range("M17").Value = 14.65250642343015
$ws.Range("N17").Value = 19.39371780311629

$ws.Range("B18").Value = 10.89228024967039
$ws.Range("C18").Value = 5.187088135132573
$ws.Range("E18").Value = 12.86783885346954
$ws.Range("F18").Value = 19.95656407809801
$ws.Range("G18").Value = 35.32865391755091
$ws.Range("H18").Value = 15.71470031420436
$ws.Range("K18").Value = 9.963094687037367
$ws.Range("M18").Value = 14.60648931849939
$ws.Range("N18").Value = 19.40957319583174

$ws.Range("B19").Value = 10.87020229729977
$ws.Range("C19").Value = 5.172919339216652
$ws.Range("E19").Value = 12.84663149887837
$ws.Range("F19").Value = 19.87204792380568
$ws.Range("G19").Value = 35.32216767082205
$ws.Range("H19").Value = 15.71746268595253
$ws.Range("K19").Value = 9.948247259359604
$ws.Range("M19").Value = 14.59093680397632
$ws.Range("N19").Value = 19.41497363289342

$ws.Range("B20").Value = 10.96936029724572
$ws.Range("C20").Value = 5.236320495873938
$ws.Range("E20").Value = 12.94205176328637
$ws.Range("F20").Value = 20.24955283636154
$ws.Range("G20").Value = 35.35212246690293
$ws.Range("H20").Value = 15.70522440047402
$ws.Range("K20").Value = 10.01502957834589
$ws.Range("M20").Value = 14.66103614059413
$ws.Range("N20").Value = 19.39079855249343

$ws.Range("B21").Value = 11.29876448611967
$ws.Range("C21").Value = 5.442959623254967
$ws.Range("E21").Value = 13.26205421336453
$ws.Range("F21").Value = 21.46857628470577
$ws.Range("G21").Value = 35.46613881198352
$ws.Range("H21").Value = 15.66748312563777
$ws.Range("K21").Value = 10.2386178340621
$ws.Range("M21").Value = 14.89829877612397
$ws.Range("N21").Value = 19.31179977320233

$ws.Range("B22").Value = 11.51086255814577
$ws.Range("C22").Value = 5.57316650679039
$ws.Range("E22").Value = 13.47035956267808
$ws.Range("F22").Value = 22.22866616901552
$ws.Range("G22").Value = 35.55057457013211
$ws.Range("H22").Value = 15.64534097715174
$ws.Range("K22").Value = 10.38389855503379
$ws.Range("M22").Value = 15.05440718873188
$ws.Range("N22").Value = 19.26182286474601

$ws.Range("B23").Value = 11.3979774163096
$ws.Range("C23").Value = 5.504121095882172
$ws.Range("E23").Value = 13.35928391051332
$ws.Range("F23").Value = 21.82633154458858
$ws.Range("G23").Value = 35.50460545354273
$ws.Range("H23").Value = 15.65692730666168
$ws.Range("K23").Value = 10.30645317298086
$ws.Range("M23").Value = 14.97100980942724
$ws.Range("N23").Value = 19.28834557147714

$ws.Range("B24").Value = 10.9639281005279
$ws.Range("C24").Value = 5.232862565091159
$ws.Range("E24").Value = 12.93681296507043
$ws.Range("F24").Value = 20.22900810905287
$ws.Range("G24").Value = 35.35042708373335
$ws.Range("H24").Value = 15.70588377858394
$ws.Range("K24").Value = 10.01136451832789
$ws.Range("M24").Value = 14.65717943464068
$ws.Range("N24").Value = 19.39211774241858

$ws.Range("B25").Value = 10.48465376942992
$ws.Range("C25").Value = 4.920080968600196
$ws.Range("E25").Value = 12.48009587335307
$ws.Range("F25").Value = 18.34778573295695
$ws.Range("G25").Value = 35.22697752144919
$ws.Range("H25").Value = 15.76948463515959
$ws.Range("K25").Value = 9.691121101945583
$ws.Range("M25").Value = 14.32477960262114
$ws.Range("N25").Value = 19.51118968789913
